$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 70 new 2-hour candle rows (r=1170..1239) after the existing data
# (which ended at row 1169). The sheet's used range/dimension grows from
# A1:F1169 to A1:F1239 automatically as a result.

# Copy the date/time style (numFmt "YYYY-MM-DD HH:MM:SS") used by column A
# from the last existing data row (A1169) so the new A-column cells match it.
$ws.Range("A1169").Copy() | Out-Null

$data = @(
    @(1170, 45534.58333333334, 0.5619, 0.5631, 0.5462, 0.5485, 55901766),
    @(1171, 45534.66666666666, 0.5485, 0.5590000000000001, 0.5463, 0.5558999999999999, 33972804),
    @(1172, 45534.75, 0.5558999999999999, 0.5629, 0.5558999999999999, 0.5618, 13458940),
    @(1173, 45534.83333333334, 0.5617, 0.5629999999999999, 0.5603, 0.5615, 2597992),
    @(1174, 45534.91666666666, 0.5658, 0.5677, 0.5657, 0.5667, 6102747),
    @(1175, 45535, 0.5667, 0.5682, 0.5653, 0.5664, 4516561),
    @(1176, 45535.08333333334, 0.5665, 0.5668, 0.5647, 0.5648, 4918291),
    @(1177, 45535.16666666666, 0.5648, 0.5687, 0.5645, 0.5671, 4486197),
    @(1178, 45535.25, 0.5669999999999999, 0.5682, 0.5664, 0.5665, 3712342),
    @(1179, 45535.33333333334, 0.5665, 0.5679999999999999, 0.5625, 0.5678, 9246523),
    @(1180, 45535.41666666666, 0.5678, 0.5721000000000001, 0.5677, 0.5692, 8427858),
    @(1181, 45535.5, 0.5693, 0.57, 0.5658, 0.5659999999999999, 6520571),
    @(1182, 45535.58333333334, 0.5659999999999999, 0.5673, 0.5644, 0.5664, 6376268),
    @(1183, 45535.66666666666, 0.5664, 0.5664, 0.5624, 0.5649, 6268903),
    @(1184, 45535.75, 0.5649, 0.5681, 0.5644, 0.5678, 4849669),
    @(1185, 45535.83333333334, 0.5678, 0.5681, 0.5661, 0.5671, 2299466),
    @(1186, 45535.91666666666, 0.5662, 0.5686, 0.5655, 0.5662, 2819279),
    @(1187, 45536, 0.5662, 0.5667, 0.5613, 0.5615, 4381696),
    @(1188, 45536.08333333334, 0.5616, 0.5618, 0.5590000000000001, 0.5608, 7543338),
    @(1189, 45536.16666666666, 0.5607, 0.5607, 0.553, 0.5545, 9564275),
    @(1190, 45536.25, 0.5545, 0.5596, 0.5541, 0.5591, 6642145),
    @(1191, 45536.33333333334, 0.5590000000000001, 0.5603, 0.5580000000000001, 0.5591, 4661558),
    @(1192, 45536.41666666666, 0.5593, 0.5593, 0.5555, 0.5580000000000001, 6459433),
    @(1193, 45536.5, 0.5581, 0.5589, 0.5528999999999999, 0.5566, 8622611),
    @(1194, 45536.58333333334, 0.5566, 0.5575, 0.5515, 0.5565, 13472172),
    @(1195, 45536.66666666666, 0.5566, 0.5583, 0.5558, 0.5558, 2045733),
    @(1196, 45536.75, 0.5551, 0.5608, 0.555, 0.5595, 9257032),
    @(1197, 45536.83333333334, 0.5595, 0.5619, 0.5553, 0.5581, 5965927),
    @(1198, 45536.91666666666, 0.5581, 0.5581, 0.5434, 0.5472, 36213741),
    @(1199, 45537, 0.5471, 0.5508, 0.5463, 0.5488, 21050608),
    @(1200, 45537.08333333334, 0.5488, 0.5503, 0.5466, 0.5501, 9542693),
    @(1201, 45537.16666666666, 0.5501, 0.5511, 0.5447, 0.5498, 10994298),
    @(1202, 45537.25, 0.5498, 0.5508999999999999, 0.5443, 0.5463, 23738519),
    @(1203, 45537.33333333334, 0.5463, 0.5543, 0.5455, 0.5535, 16183905),
    @(1204, 45537.41666666666, 0.5534, 0.5576, 0.5521, 0.5553, 21525558),
    @(1205, 45537.5, 0.5548999999999999, 0.5561, 0.5526, 0.5528, 12566402),
    @(1206, 45537.58333333334, 0.5528, 0.5531, 0.5528, 0.5528999999999999, 123223),
    @(1207, 45537.66666666666, 0.5582, 0.5609, 0.5555, 0.5608, 8942406),
    @(1208, 45537.75, 0.5608, 0.5679999999999999, 0.5596, 0.5679, 13478318),
    @(1209, 45537.83333333334, 0.5678, 0.5688, 0.5644, 0.5665, 12759624),
    @(1210, 45537.91666666666, 0.5664, 0.5685, 0.5658, 0.5677, 7264965),
    @(1211, 45538, 0.5677, 0.569, 0.5656, 0.569, 10663816),
    @(1212, 45538.08333333334, 0.5689, 0.5724, 0.5666, 0.5667, 13470119),
    @(1213, 45538.16666666666, 0.5666, 0.5676, 0.5649, 0.5675, 9504724),
    @(1214, 45538.25, 0.5674, 0.5695, 0.5656, 0.5692, 21607393),
    @(1215, 45538.33333333334, 0.5691000000000001, 0.572, 0.5658, 0.57, 12696430),
    @(1216, 45538.41666666666, 0.57, 0.5724, 0.5679999999999999, 0.5701000000000001, 18111837),
    @(1217, 45538.5, 0.57, 0.5723, 0.5605, 0.5625, 22835905),
    @(1218, 45538.58333333334, 0.5624, 0.5649999999999999, 0.5571, 0.5593, 27717123),
    @(1219, 45538.66666666666, 0.5592, 0.5648, 0.5586, 0.5617, 15211846),
    @(1220, 45538.75, 0.5617, 0.5659, 0.5607, 0.5659, 10119440),
    @(1221, 45538.83333333334, 0.5659, 0.5677, 0.5648, 0.5661, 7987449),
    @(1222, 45538.91666666666, 0.5659999999999999, 0.5662, 0.5565, 0.5572, 9855401),
    @(1223, 45539, 0.5572, 0.5604, 0.5326, 0.5501, 82012366),
    @(1224, 45539.08333333334, 0.5502, 0.5548999999999999, 0.5489000000000001, 0.553, 13780595),
    @(1225, 45539.16666666666, 0.553, 0.5538, 0.5497, 0.55, 9914187),
    @(1226, 45539.25, 0.55, 0.5556, 0.5499000000000001, 0.5545, 23636353),
    @(1227, 45539.33333333334, 0.5546, 0.5570000000000001, 0.5528999999999999, 0.553, 17096355),
    @(1228, 45539.41666666666, 0.5531, 0.5533, 0.549, 0.5505, 26169811),
    @(1229, 45539.5, 0.5505, 0.5543, 0.5472, 0.5537, 24596822),
    @(1230, 45539.58333333334, 0.5536, 0.5629999999999999, 0.5516, 0.5621, 33485074),
    @(1231, 45539.66666666666, 0.5621, 0.5648, 0.5585, 0.5590000000000001, 16210340),
    @(1232, 45539.75, 0.5590000000000001, 0.5605, 0.5564, 0.5567, 6370880),
    @(1233, 45539.83333333334, 0.5578, 0.5580000000000001, 0.5574, 0.5577, 118399),
    @(1234, 45539.91666666666, 0.5586, 0.5612, 0.5572, 0.5579, 10296007),
    @(1235, 45540, 0.5579, 0.5598, 0.5572, 0.5585, 14721682),
    @(1236, 45540.08333333334, 0.5585, 0.5589, 0.5517, 0.5528, 14556876),
    @(1237, 45540.16666666666, 0.5528999999999999, 0.5542, 0.5523, 0.553, 6786478),
    @(1238, 45540.25, 0.553, 0.5547, 0.5499000000000001, 0.5535, 29195527),
    @(1239, 45540.33333333334, 0.5535, 0.5547, 0.5498, 0.5499000000000001, 6963768)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("A" + $r).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
}